$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "END OF PAID HOURS" -> "END OF 2nd SET HOURS" (shared string used by D23)
$ws.Range("D23").Value = "END OF 2nd SET HOURS"

# 2. Widen column D (bestFit width 15.33 -> 21)
$ws.Columns.Item(4).ColumnWidth = 20.17

# 3. Add a new logged session in row 24: lesson #22, date 2/15/2018, 2 hours
$ws.Range("B23").Copy($ws.Range("B24"))
$ws.Range("C23").Copy($ws.Range("C24"))
$ws.Range("A24").Value = 22
$ws.Range("B24").Value2 = 43146
$ws.Range("C24").Value = 2
